# This script reorders / trims the employee table in the "Worksheet" sheet.
# Plan:
#  1. Copy the 11 surviving data rows (everything except the "Samsul Huda" row)
#     from their current positions into a scratch area far below the table,
#     already arranged in the desired final order. Using a Range.Copy keeps
#     the original cell types (shared-string text, numbers, ...) intact, so
#     text such as phone numbers with a leading zero is preserved exactly.
#  2. Clear out the old table rows (2-13).
#  3. Copy the scratch rows back into the table body (rows 2-12).
#  4. Clear the scratch rows.
#  5. Patch the handful of numeric cells whose value actually changed
#     (ID column for Tolkha Hasan / Wawan Dwi Prasetyo, ID Kota for
#     Achmad Chadil Auwfar).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source row (in the original layout) for each row of the final table, in order.
$sourceRows = @(3, 4, 5, 7, 8, 9, 11, 12, 13, 6, 10)

$scratchStart = 200
for ($i = 0; $i -lt $sourceRows.Length; $i++) {
    $src = $sourceRows[$i]
    $dst = $scratchStart + $i
    $ws.Range("A$src`:G$src").Copy($ws.Range("A$dst`:G$dst"))
}

# Wipe the old table body (header row 1 is left untouched).
$ws.Range("A2:G13").ClearContents()

# Move the freshly ordered rows back into the table body.
for ($i = 0; $i -lt $sourceRows.Length; $i++) {
    $src = $scratchStart + $i
    $dst = 2 + $i
    $ws.Range("A$src`:G$src").Copy($ws.Range("A$dst`:G$dst"))
}

# Remove the scratch data.
$ws.Range("A$scratchStart`:G$($scratchStart + $sourceRows.Length - 1)").ClearContents()

# Fix the few cells whose value genuinely changed versus the source row.
# Row 8  -> Achmad Chadil Auwfar : ID Kota 1 -> 2
$ws.Cells.Item(8, 4).Value = 2
# Row 11 -> Tolkha Hasan         : ID 2 -> 1
$ws.Cells.Item(11, 1).Value = 1
# Row 12 -> Wawan Dwi Prasetyo   : ID 6 -> 2
$ws.Cells.Item(12, 1).Value = 2
